# The "Carport_Liste" sheet header in E1 previously read "Formel / Info";
# update it to just "Formel" (matching the shared-string change in the diff).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Carport_Liste")

$ws.Range("E1").Value = "Formel"

# Reflect the author's selection landing on the edited cell.
$ws.Range("E1").Select()
